$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (RF)
$ws.Range("B3").Value = 0.825
$ws.Range("C3").Value = 0.904
$ws.Range("D3").Value = 0.646
$ws.Range("E3").Value = 0.909
$ws.Range("F3").Value = 0.901
$ws.Range("G3").Value = 0.099
$ws.Range("H3").Value = 0.315
$ws.Range("I3").Value = 0.234
$ws.Range("J3").Value = 0.974

# Row 4 (NN)
$ws.Range("E4").Value = 0.732
$ws.Range("F4").Value = 0.71
$ws.Range("G4").Value = 0.292
$ws.Range("H4").Value = 0.54
$ws.Range("I4").Value = 0.412
$ws.Range("J4").Value = 0.915

# Row 5 (RNN)
$ws.Range("E5").Value = 0.552
$ws.Range("F5").Value = 0.533
$ws.Range("G5").Value = 0.489
$ws.Range("H5").Value = 0.699
$ws.Range("I5").Value = 0.545
$ws.Range("J5").Value = 0.793

# Row 6 (Ensemble)
$ws.Range("E6").Value = 0.591
$ws.Range("F6").Value = 0.5570000000000001
$ws.Range("G6").Value = 0.446
$ws.Range("H6").Value = 0.668
$ws.Range("I6").Value = 0.491
$ws.Range("J6").Value = 0.868

$wb.Save()
